# Applies the price/volume refresh from the Tue Aug 27 05:52:11 UTC 2024
# GitHub Actions run, including the Mantle/Hedera row re-ranking (rows 45-46).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.076.20"
$ws.Range("E2").Value = "  -1.37%  "
$ws.Range("D3").Value = "2.689.64"
$ws.Range("E3").Value = "  -1.82%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'557.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.71%  "
$ws.Range("D6").Value = "'159.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -0.79%  "
$ws.Range("E9").Value = "  -2.66%  "
$ws.Range("E10").Value = "  -1.87%  "
$ws.Range("D11").Value = "'0.372"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.09%  "
$ws.Range("E12").Value = "  -5.45%  "
$ws.Range("D13").Value = "3.165.33"
$ws.Range("E13").Value = "  -1.83%  "
$ws.Range("E14").Value = "  -1.54%  "
$ws.Range("D15").Value = "62.979.28"
$ws.Range("E15").Value = "  -1.26%  "
$ws.Range("E16").Value = "  -1.45%  "
$ws.Range("D17").Value = "2.690.23"
$ws.Range("E17").Value = "  -1.92%  "
$ws.Range("D18").Value = "'12.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'4.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.23%  "
$ws.Range("D20").Value = "'346.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.91%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("E23").Value = "  -1.68%  "
$ws.Range("D24").Value = "'63.56"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.13%  "
$ws.Range("E25").Value = "  -0.68%  "
$ws.Range("D26").Value = "'0.997"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("E27").Value = "  -1.76%  "
$ws.Range("E28").Value = "  +9.20%  "
$ws.Range("E29").Value = "  -4.80%  "
$ws.Range("D30").Value = "'7.29"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.84%  "
$ws.Range("E31").Value = "  -0.12%  "
$ws.Range("D32").Value = "'164.54"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.32%  "
$ws.Range("D33").Value = "'4.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.09%  "
$ws.Range("E34").Value = "  +1.45%  "
$ws.Range("D35").Value = "'0.998"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  -2.59%  "
$ws.Range("E37").Value = "  +0.22%  "
$ws.Range("D38").Value = "'361.45"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.66%  "
$ws.Range("D39").Value = "'6.48"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.16%  "
$ws.Range("E40").Value = "  -2.21%  "
$ws.Range("E41").Value = "  -1.32%  "
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("D43").Value = "'21.23"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.81%  "
$ws.Range("D44").Value = "'20.46"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.95%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").Value = "'0.0566"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.95%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "'0.621"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.35%  "
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("D49").Value = "'0.0245"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.22%  "
$ws.Range("D50").Value = "'130.07"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.19%  "
$ws.Range("E51").Value = "  -2.78%  "
